# Update the "dSF" column (F) values for several rows, per repull/recalculation
# of the underlying data (commit: "repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F4"  = -4
    "F8"  = 9
    "F9"  = -6
    "F11" = 1
    "F13" = 1
    "F14" = -5
    "F17" = -4
    "F21" = -3
    "F23" = 0
    "F26" = 4
    "F31" = 0
    "F34" = -5
    "F38" = -1
    "F39" = 1
    "F41" = -3
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
